# Add the missing student/trainee names into column A (rows 2-10) of the
# "Sagar Lab" sheet. These rows already held the user/password/URL columns
# (B-D); the author filled in the corresponding "Name" values in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ordered list (row 2 .. row 10) so shared-string insertion order matches
# the order the names were typed in.
$names = @(
    "Gautham G",
    "Narendra kumar",
    "Nageswara Rao",
    "Gautham T",
    "Sunanda",
    "Krishna Kumari",
    "Girija Shankar",
    "Shivam Singh",
    "Nirmal"
)

$startRow = 2
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $names[$i]
}

# Reflect the final cursor/selection position left by the editing session.
$ws.Range("A10").Select()
